$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Formula = "=`"39.962.53`""
$ws.Range("D2").Copy($ws.Range("D2")) | Out-Null
$ws.Range("D2").PasteSpecial(-4163) | Out-Null
$ws.Range("E2").Value = "  +1.32%  "

$ws.Range("D3").Formula = "=`"2.205.86`""
$ws.Range("D3").Copy($ws.Range("D3")) | Out-Null
$ws.Range("D3").PasteSpecial(-4163) | Out-Null
$ws.Range("E3").Value = "  +2.04%  "

$ws.Range("E4").Value = "  +0.00%  "

$ws.Range("D5").Formula = "=`"228.30`""
$ws.Range("D5").Copy($ws.Range("D5")) | Out-Null
$ws.Range("D5").PasteSpecial(-4163) | Out-Null
$ws.Range("E5").Value = "  -0.30%  "

$ws.Range("D6").Formula = "=`"0.629`""
$ws.Range("D6").Copy($ws.Range("D6")) | Out-Null
$ws.Range("D6").PasteSpecial(-4163) | Out-Null
$ws.Range("E6").Value = "  +0.96%  "

$ws.Range("D7").Formula = "=`"63.71`""
$ws.Range("D7").Copy($ws.Range("D7")) | Out-Null
$ws.Range("D7").PasteSpecial(-4163) | Out-Null
$ws.Range("E7").Value = "  +0.74%  "

$ws.Range("E8").Value = "  +0.05%  "

$ws.Range("D9").Formula = "=`"0.395`""
$ws.Range("D9").Copy($ws.Range("D9")) | Out-Null
$ws.Range("D9").PasteSpecial(-4163) | Out-Null
$ws.Range("E9").Value = "  -0.16%  "

$ws.Range("D10").Formula = "=`"0.0856`""
$ws.Range("D10").Copy($ws.Range("D10")) | Out-Null
$ws.Range("D10").PasteSpecial(-4163) | Out-Null
$ws.Range("E10").Value = "  -0.63%  "

$ws.Range("E11").Value = "  +0.71%  "

$ws.Range("D12").Formula = "=`"16.12`""
$ws.Range("D12").Copy($ws.Range("D12")) | Out-Null
$ws.Range("D12").PasteSpecial(-4163) | Out-Null
$ws.Range("E12").Value = "  +0.27%  "

$ws.Range("D13").Formula = "=`"2.530.55`""
$ws.Range("D13").Copy($ws.Range("D13")) | Out-Null
$ws.Range("D13").PasteSpecial(-4163) | Out-Null
$ws.Range("E13").Value = "  +2.00%  "

$ws.Range("D14").Formula = "=`"22.13`""
$ws.Range("D14").Copy($ws.Range("D14")) | Out-Null
$ws.Range("D14").PasteSpecial(-4163) | Out-Null
$ws.Range("E14").Value = "  -0.51%  "

$ws.Range("D15").Formula = "=`"0.822`""
$ws.Range("D15").Copy($ws.Range("D15")) | Out-Null
$ws.Range("D15").PasteSpecial(-4163) | Out-Null
$ws.Range("E15").Value = "  +0.51%  "

$ws.Range("D16").Formula = "=`"5.60`""
$ws.Range("D16").Copy($ws.Range("D16")) | Out-Null
$ws.Range("D16").PasteSpecial(-4163) | Out-Null
$ws.Range("E16").Value = "  +0.38%  "

$ws.Range("D17").Formula = "=`"2.194.55`""
$ws.Range("D17").Copy($ws.Range("D17")) | Out-Null
$ws.Range("D17").PasteSpecial(-4163) | Out-Null
$ws.Range("E17").Value = "  +2.28%  "

$ws.Range("D18").Formula = "=`"39.912.58`""
$ws.Range("D18").Copy($ws.Range("D18")) | Out-Null
$ws.Range("D18").PasteSpecial(-4163) | Out-Null
$ws.Range("E18").Value = "  +1.34%  "

$ws.Range("D19").Formula = "=`"0.0₃0916`""
$ws.Range("D19").Copy($ws.Range("D19")) | Out-Null
$ws.Range("D19").PasteSpecial(-4163) | Out-Null
$ws.Range("E19").Value = "  +7.52%  "

$ws.Range("D20").Formula = "=`"72.17`""
$ws.Range("D20").Copy($ws.Range("D20")) | Out-Null
$ws.Range("D20").PasteSpecial(-4163) | Out-Null
$ws.Range("E20").Value = "  -0.25%  "

$ws.Range("D21").Formula = "=`"6.08`""
$ws.Range("D21").Copy($ws.Range("D21")) | Out-Null
$ws.Range("D21").PasteSpecial(-4163) | Out-Null
$ws.Range("E21").Value = "  -1.12%  "

$ws.Range("D22").Formula = "=`"231.59`""
$ws.Range("D22").Copy($ws.Range("D22")) | Out-Null
$ws.Range("D22").PasteSpecial(-4163) | Out-Null
$ws.Range("E22").Value = "  +1.39%  "

$ws.Range("E23").Value = "  +0.02%  "

$ws.Range("E24").Value = "  +0.00%  "

$ws.Range("D25").Formula = "=`"2.36`""
$ws.Range("D25").Copy($ws.Range("D25")) | Out-Null
$ws.Range("D25").PasteSpecial(-4163) | Out-Null
$ws.Range("E25").Value = "  -0.10%  "

$ws.Range("B26").Value = "Cosmos"
$ws.Range("C26").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D26").Formula = "=`"9.67`""
$ws.Range("D26").Copy($ws.Range("D26")) | Out-Null
$ws.Range("D26").PasteSpecial(-4163) | Out-Null
$ws.Range("E26").Value = "  +0.10%  "

$ws.Range("B27").Value = "Monero"
$ws.Range("C27").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D27").Formula = "=`"171.35`""
$ws.Range("D27").Copy($ws.Range("D27")) | Out-Null
$ws.Range("D27").PasteSpecial(-4163) | Out-Null
$ws.Range("E27").Value = "  -0.51%  "

$ws.Range("E28").Value = "  +1.25%  "

$ws.Range("D29").Formula = "=`"1.46`""
$ws.Range("D29").Copy($ws.Range("D29")) | Out-Null
$ws.Range("D29").PasteSpecial(-4163) | Out-Null
$ws.Range("E29").Value = "  +3.00%  "

$ws.Range("D30").Formula = "=`"20.12`""
$ws.Range("D30").Copy($ws.Range("D30")) | Out-Null
$ws.Range("D30").PasteSpecial(-4163) | Out-Null
$ws.Range("E30").Value = "  +2.37%  "

$ws.Range("D31").Formula = "=`"2.71`""
$ws.Range("D31").Copy($ws.Range("D31")) | Out-Null
$ws.Range("D31").PasteSpecial(-4163) | Out-Null
$ws.Range("E31").Value = "  +5.21%  "

$ws.Range("E32").Value = "  +0.94%  "

$ws.Range("D33").Formula = "=`"4.58`""
$ws.Range("D33").Copy($ws.Range("D33")) | Out-Null
$ws.Range("D33").PasteSpecial(-4163) | Out-Null
$ws.Range("E33").Value = "  -1.85%  "

$ws.Range("D34").Formula = "=`"4.72`""
$ws.Range("D34").Copy($ws.Range("D34")) | Out-Null
$ws.Range("D34").PasteSpecial(-4163) | Out-Null
$ws.Range("E34").Value = "  -1.95%  "

$ws.Range("D35").Formula = "=`"7.04`""
$ws.Range("D35").Copy($ws.Range("D35")) | Out-Null
$ws.Range("D35").PasteSpecial(-4163) | Out-Null
$ws.Range("E35").Value = "  +0.07%  "

$ws.Range("D36").Formula = "=`"0.0624`""
$ws.Range("D36").Copy($ws.Range("D36")) | Out-Null
$ws.Range("D36").PasteSpecial(-4163) | Out-Null
$ws.Range("E36").Value = "  +0.27%  "

$ws.Range("D37").Formula = "=`"3.86`""
$ws.Range("D37").Copy($ws.Range("D37")) | Out-Null
$ws.Range("D37").PasteSpecial(-4163) | Out-Null
$ws.Range("E37").Value = "  +8.97%  "

$ws.Range("D38").Formula = "=`"2.46`""
$ws.Range("D38").Copy($ws.Range("D38")) | Out-Null
$ws.Range("D38").PasteSpecial(-4163) | Out-Null
$ws.Range("E38").Value = "  +1.03%  "

$ws.Range("D39").Formula = "=`"0.998`""
$ws.Range("D39").Copy($ws.Range("D39")) | Out-Null
$ws.Range("D39").PasteSpecial(-4163) | Out-Null
$ws.Range("E39").Value = "  -0.10%  "

$ws.Range("D40").Formula = "=`"5.01`""
$ws.Range("D40").Copy($ws.Range("D40")) | Out-Null
$ws.Range("D40").PasteSpecial(-4163) | Out-Null
$ws.Range("E40").Value = "  +18.38%  "

$ws.Range("D41").Formula = "=`"103.65`""
$ws.Range("D41").Copy($ws.Range("D41")) | Out-Null
$ws.Range("D41").PasteSpecial(-4163) | Out-Null
$ws.Range("E41").Value = "  +0.02%  "

$ws.Range("D42").Formula = "=`"0.0230`""
$ws.Range("D42").Copy($ws.Range("D42")) | Out-Null
$ws.Range("D42").PasteSpecial(-4163) | Out-Null
$ws.Range("E42").Value = "  -0.48%  "

$ws.Range("D43").Formula = "=`"17.96`""
$ws.Range("D43").Copy($ws.Range("D43")) | Out-Null
$ws.Range("D43").PasteSpecial(-4163) | Out-Null
$ws.Range("E43").Value = "  -1.22%  "

$ws.Range("B44").Value = "TrustWalletToken"
$ws.Range("C44").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D44").Formula = "=`"1.23`""
$ws.Range("D44").Copy($ws.Range("D44")) | Out-Null
$ws.Range("D44").PasteSpecial(-4163) | Out-Null
$ws.Range("E44").Value = "  +3.38%  "

$ws.Range("B45").Value = "Maker"
$ws.Range("C45").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D45").Formula = "=`"1.520.08`""
$ws.Range("D45").Copy($ws.Range("D45")) | Out-Null
$ws.Range("D45").PasteSpecial(-4163) | Out-Null
$ws.Range("E45").Value = "  -0.65%  "

$ws.Range("D46").Formula = "=`"8.04`""
$ws.Range("D46").Copy($ws.Range("D46")) | Out-Null
$ws.Range("D46").PasteSpecial(-4163) | Out-Null
$ws.Range("E46").Value = "  +3.46%  "

$ws.Range("B47").Value = "ARBITRUM"
$ws.Range("C47").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D47").Formula = "=`"1.11`""
$ws.Range("D47").Copy($ws.Range("D47")) | Out-Null
$ws.Range("D47").PasteSpecial(-4163) | Out-Null
$ws.Range("E47").Value = "  +0.21%  "

$ws.Range("B48").Value = "Cronos"
$ws.Range("C48").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D48").Formula = "=`"0.0926`""
$ws.Range("D48").Copy($ws.Range("D48")) | Out-Null
$ws.Range("D48").PasteSpecial(-4163) | Out-Null
$ws.Range("E48").Value = "  -0.68%  "

$ws.Range("B49").Value = "HuobiToken"
$ws.Range("C49").Value = "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
$ws.Range("D49").Formula = "=`"2.81`""
$ws.Range("D49").Copy($ws.Range("D49")) | Out-Null
$ws.Range("D49").PasteSpecial(-4163) | Out-Null
$ws.Range("E49").Value = "  -0.06%  "

$ws.Range("D50").Formula = "=`"0.000194`""
$ws.Range("D50").Copy($ws.Range("D50")) | Out-Null
$ws.Range("D50").PasteSpecial(-4163) | Out-Null
$ws.Range("E50").Value = "  +32.83%  "

$ws.Range("D51").Formula = "=`"2.409.33`""
$ws.Range("D51").Copy($ws.Range("D51")) | Out-Null
$ws.Range("D51").PasteSpecial(-4163) | Out-Null
$ws.Range("E51").Value = "  +1.91%  "
